# Update the "Förändrad" (changed) date for rows 2-11, column C,
# from 2023-10-09 (45208) to 2023-10-13 (45212).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C11").Value = 45212
